$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 2
Set-TextValue $ws.Range('D2') '42.813.98'
Set-TextValue $ws.Range('E2') '  +3.67%  '

# Row 3
Set-TextValue $ws.Range('D3') '2.256.96'
Set-TextValue $ws.Range('E3') '  +3.35%  '

# Row 4
Set-TextValue $ws.Range('E4') '  +0.02%  '

# Row 5
Set-TextValue $ws.Range('D5') '253.54'
Set-TextValue $ws.Range('E5') '  -0.72%  '

# Row 6
Set-TextValue $ws.Range('D6') '0.636'
Set-TextValue $ws.Range('E6') '  +1.30%  '

# Row 7
Set-TextValue $ws.Range('D7') '70.88'
Set-TextValue $ws.Range('E7') '  +3.75%  '

# Row 8
Set-TextValue $ws.Range('E8') '  -0.12%  '

# Row 9
Set-TextValue $ws.Range('D9') '0.651'
Set-TextValue $ws.Range('E9') '  +12.28%  '

# Row 10
Set-TextValue $ws.Range('D10') '41.53'
Set-TextValue $ws.Range('E10') '  +8.67%  '

# Row 11
Set-TextValue $ws.Range('D11') '59.54'
Set-TextValue $ws.Range('E11') '  +1.06%  '

# Row 12
Set-TextValue $ws.Range('D12') '0.0961'
Set-TextValue $ws.Range('E12') '  +2.53%  '

# Row 13
Set-TextValue $ws.Range('E13') '  +2.92%  '

# Row 14
Set-TextValue $ws.Range('D14') '0.105'
Set-TextValue $ws.Range('E14') '  +0.35%  '

# Row 15
Set-TextValue $ws.Range('D15') '2.592.89'
Set-TextValue $ws.Range('E15') '  +3.41%  '

# Row 16
Set-TextValue $ws.Range('B16') 'Chainlink'
Set-TextValue $ws.Range('C16') 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range('D16') '14.88'
Set-TextValue $ws.Range('E16') '  +2.45%  '

# Row 17
Set-TextValue $ws.Range('B17') 'Polygon'
Set-TextValue $ws.Range('C17') 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue $ws.Range('D17') '0.889'
Set-TextValue $ws.Range('E17') '  +1.24%  '

# Row 18
Set-TextValue $ws.Range('D18') '2.252.87'
Set-TextValue $ws.Range('E18') '  +3.16%  '

# Row 19
Set-TextValue $ws.Range('D19') '42.730.86'
Set-TextValue $ws.Range('E19') '  +3.65%  '

# Row 20
Set-TextValue $ws.Range('D20') '0.0₃0984'
Set-TextValue $ws.Range('E20') '  +2.95%  '

# Row 21
Set-TextValue $ws.Range('D21') '6.27'
Set-TextValue $ws.Range('E21') '  +1.32%  '

# Row 22
Set-TextValue $ws.Range('D22') '73.12'
Set-TextValue $ws.Range('E22') '  +1.62%  '

# Row 23
Set-TextValue $ws.Range('D23') '236.26'
Set-TextValue $ws.Range('E23') '  +1.59%  '

# Row 24
Set-TextValue $ws.Range('D24') '2.14'
Set-TextValue $ws.Range('E24') '  +4.56%  '

# Row 25
Set-TextValue $ws.Range('D25') '3.99'
Set-TextValue $ws.Range('E25') '  +1.08%  '

# Row 26
Set-TextValue $ws.Range('D26') '11.73'
Set-TextValue $ws.Range('E26') '  -1.79%  '

# Row 27
Set-TextValue $ws.Range('E27') '  +0.04%  '

# Row 28
Set-TextValue $ws.Range('D28') '2.45'
Set-TextValue $ws.Range('E28') '  -3.69%  '

# Row 29
Set-TextValue $ws.Range('E29') '  -1.73%  '

# Row 30
Set-TextValue $ws.Range('E30') '  +1.62%  '

# Row 31
Set-TextValue $ws.Range('D31') '167.78'
Set-TextValue $ws.Range('E31') '  -0.83%  '

# Row 32
Set-TextValue $ws.Range('D32') '21.04'
Set-TextValue $ws.Range('E32') '  +1.78%  '

# Row 33
Set-TextValue $ws.Range('B33') 'Kaspa'
Set-TextValue $ws.Range('C33') 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range('D33') '0.129'
Set-TextValue $ws.Range('E33') '  +10.18%  '

# Row 34
Set-TextValue $ws.Range('B34') 'InternetComputer(DFINITY)'
Set-TextValue $ws.Range('C34') 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextValue $ws.Range('D34') '6.16'
Set-TextValue $ws.Range('E34') '  +12.74%  '

# Row 35
Set-TextValue $ws.Range('D35') '0.0790'
Set-TextValue $ws.Range('E35') '  +5.85%  '

# Row 36
Set-TextValue $ws.Range('E36') '  +0.97%  '

# Row 37
Set-TextValue $ws.Range('D37') '28.04'
Set-TextValue $ws.Range('E37') '  +2.79%  '

# Row 38
Set-TextValue $ws.Range('D38') '4.72'

# Row 39
Set-TextValue $ws.Range('D39') '4.17'
Set-TextValue $ws.Range('E39') '  -0.81%  '

# Row 40
Set-TextValue $ws.Range('D40') '0.0318'
Set-TextValue $ws.Range('E40') '  +6.14%  '

# Row 41
Set-TextValue $ws.Range('E41') '  +3.41%  '

# Row 42
Set-TextValue $ws.Range('D42') '12.48'
Set-TextValue $ws.Range('E42') '  -1.23%  '

# Row 43
Set-TextValue $ws.Range('D43') '5.82'
Set-TextValue $ws.Range('E43') '  +2.17%  '

# Row 44
Set-TextValue $ws.Range('D44') '64.66'
Set-TextValue $ws.Range('E44') '  -0.16%  '

# Row 45
Set-TextValue $ws.Range('B45') 'FTXToken'
Set-TextValue $ws.Range('C45') 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
Set-TextValue $ws.Range('D45') '5.00'
Set-TextValue $ws.Range('E45') '  -1.81%  '

# Row 46
Set-TextValue $ws.Range('B46') 'Algorand'
Set-TextValue $ws.Range('C46') 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
Set-TextValue $ws.Range('D46') '0.204'
Set-TextValue $ws.Range('E46') '  +1.37%  '

# Row 47
Set-TextValue $ws.Range('D47') '8.93'
Set-TextValue $ws.Range('E47') '  +3.08%  '

# Row 48
Set-TextValue $ws.Range('D48') '0.103'
Set-TextValue $ws.Range('E48') '  +1.14%  '

# Row 49
Set-TextValue $ws.Range('E49') '  +5.14%  '

# Row 50
Set-TextValue $ws.Range('D50') '0.998'
Set-TextValue $ws.Range('E50') '  -0.47%  '

# Row 51
Set-TextValue $ws.Range('D51') '4.47'
Set-TextValue $ws.Range('E51') '  +3.63%  '
